# The benchmark table was reshuffled: row 1 holds the same 12 dataset headers
# but in a new column order, rows 2-6 hold the same 5 classic models (reordered;
# ComplementNB re-measured after reclassing the problem to 3 classes), a new
# "setfit" row was measured, and the SVM row (still missing 3 datasets) moved
# down to row 8 to make room for "setfit".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 2).Value = "reddit"
$ws.Cells.Item(1, 3).Value = "patio_lawn_garden"
$ws.Cells.Item(1, 4).Value = "twitter"
$ws.Cells.Item(1, 5).Value = "lidl"
$ws.Cells.Item(1, 6).Value = "automotive"
$ws.Cells.Item(1, 7).Value = "ikea_reviews"
$ws.Cells.Item(1, 8).Value = "luxury_beauty"
$ws.Cells.Item(1, 9).Value = "instant_video"
$ws.Cells.Item(1, 10).Value = "musical_instruments"
$ws.Cells.Item(1, 11).Value = "office_products"
$ws.Cells.Item(1, 12).Value = "hotel"
$ws.Cells.Item(1, 13).Value = "drugs"

# Row 2
$ws.Cells.Item(2, 1).Value = "ComplementNB"
$ws.Cells.Item(2, 2).Value = 1.368598937988281
$ws.Cells.Item(2, 3).Value = 4.545320749282837
$ws.Cells.Item(2, 4).Value = 0.8519363403320312
$ws.Cells.Item(2, 5).Value = 1.694851636886597
$ws.Cells.Item(2, 6).Value = 5.123377323150635
$ws.Cells.Item(2, 7).Value = 0.6928644180297852
$ws.Cells.Item(2, 8).Value = 6.780089616775513
$ws.Cells.Item(2, 9).Value = 7.883801460266113
$ws.Cells.Item(2, 10).Value = 2.174875736236572
$ws.Cells.Item(2, 11).Value = 17.53499841690063
$ws.Cells.Item(2, 12).Value = 1.462414979934692
$ws.Cells.Item(2, 13).Value = 10.10558819770813

# Row 3
$ws.Cells.Item(3, 1).Value = "Decision Tree"
$ws.Cells.Item(3, 2).Value = 34.25000977516174
$ws.Cells.Item(3, 3).Value = 97.10490870475769
$ws.Cells.Item(3, 4).Value = 21.49746203422546
$ws.Cells.Item(3, 5).Value = 46.26797533035278
$ws.Cells.Item(3, 6).Value = 90.58149528503418
$ws.Cells.Item(3, 7).Value = 16.56105422973633
$ws.Cells.Item(3, 8).Value = 150.0110960006714
$ws.Cells.Item(3, 9).Value = 195.7523448467255
$ws.Cells.Item(3, 10).Value = 46.93049097061157
$ws.Cells.Item(3, 11).Value = 389.5512316226959
$ws.Cells.Item(3, 12).Value = 34.79415583610535
$ws.Cells.Item(3, 13).Value = 241.411954164505

# Row 4
$ws.Cells.Item(4, 1).Value = "LR"
$ws.Cells.Item(4, 2).Value = 170.3610789775848
$ws.Cells.Item(4, 3).Value = 430.759304523468
$ws.Cells.Item(4, 4).Value = 60.14369559288025
$ws.Cells.Item(4, 5).Value = 311.0737907886505
$ws.Cells.Item(4, 6).Value = 386.1109373569489
$ws.Cells.Item(4, 7).Value = 61.5999059677124
$ws.Cells.Item(4, 8).Value = 614.8274817466736
$ws.Cells.Item(4, 9).Value = 1177.187592983246
$ws.Cells.Item(4, 10).Value = 190.9322052001953
$ws.Cells.Item(4, 11).Value = 2276.936871528625
$ws.Cells.Item(4, 12).Value = 120.3800938129425
$ws.Cells.Item(4, 13).Value = 2041.393538236618

# Row 5
$ws.Cells.Item(5, 1).Value = "MultinomialNB"
$ws.Cells.Item(5, 2).Value = 1.360639572143555
$ws.Cells.Item(5, 3).Value = 4.72423529624939
$ws.Cells.Item(5, 4).Value = 0.9022314548492432
$ws.Cells.Item(5, 5).Value = 1.595834732055664
$ws.Cells.Item(5, 6).Value = 5.115803956985474
$ws.Cells.Item(5, 7).Value = 0.6823203563690186
$ws.Cells.Item(5, 8).Value = 6.495334625244141
$ws.Cells.Item(5, 9).Value = 7.96331787109375
$ws.Cells.Item(5, 10).Value = 2.265194177627563
$ws.Cells.Item(5, 11).Value = 17.22799348831177
$ws.Cells.Item(5, 12).Value = 1.443359136581421
$ws.Cells.Item(5, 13).Value = 9.757237672805786

# Row 6
$ws.Cells.Item(6, 1).Value = "RF"
$ws.Cells.Item(6, 2).Value = 197.9030184745789
$ws.Cells.Item(6, 3).Value = 319.9789986610413
$ws.Cells.Item(6, 4).Value = 125.0073599815369
$ws.Cells.Item(6, 5).Value = 264.1496708393097
$ws.Cells.Item(6, 6).Value = 293.9082410335541
$ws.Cells.Item(6, 7).Value = 122.2454059123993
$ws.Cells.Item(6, 8).Value = 611.0661752223969
$ws.Cells.Item(6, 9).Value = 730.914412021637
$ws.Cells.Item(6, 10).Value = 176.7431120872498
$ws.Cells.Item(6, 11).Value = 1198.608413696289
$ws.Cells.Item(6, 12).Value = 159.0874509811401
$ws.Cells.Item(6, 13).Value = 880.04869556427

# Row 7
$ws.Cells.Item(7, 1).Value = "setfit"
$ws.Cells.Item(7, 2).Value = 31.68623161315918
$ws.Cells.Item(7, 3).Value = 223.4237172603607
$ws.Cells.Item(7, 4).Value = 23.57607388496399
$ws.Cells.Item(7, 5).Value = 34.34165120124817
$ws.Cells.Item(7, 6).Value = 79.33475589752197
$ws.Cells.Item(7, 7).Value = 41.35797667503357
$ws.Cells.Item(7, 8).Value = 403.7233846187592
$ws.Cells.Item(7, 9).Value = 257.0964665412903
$ws.Cells.Item(7, 10).Value = 79.90693593025208
$ws.Cells.Item(7, 11).Value = 467.1691646575928
$ws.Cells.Item(7, 12).Value = 199.0017547607422
$ws.Cells.Item(7, 13).Value = 74.40118527412415

# Row 8
$ws.Cells.Item(8, 1).Value = "SVM"
$ws.Cells.Item(8, 2).Value = $null
$ws.Cells.Item(8, 3).Value = $null
$ws.Cells.Item(8, 4).Value = $null
$ws.Cells.Item(8, 5).Value = 84303.21725845337
$ws.Cells.Item(8, 6).Value = 16120.23056173325
$ws.Cells.Item(8, 7).Value = 4361.065144777298
$ws.Cells.Item(8, 8).Value = 43334.6003139019
$ws.Cells.Item(8, 9).Value = 50211.97952365875
$ws.Cells.Item(8, 10).Value = 3974.249224185944
$ws.Cells.Item(8, 11).Value = $null
$ws.Cells.Item(8, 12).Value = 7319.897396087646
$ws.Cells.Item(8, 13).Value = 228744.987988472

# Row 8 is brand new, so its row label (A8) needs the same bold, bordered,
# centered style already used by every other row label / column header.
# Copy the formatting (not the value) from an existing label cell so the
# same shared style entry is reused instead of a new one being created.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(8, 1).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
